$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (including the date number format) from the row above
# down into the new row 7 before filling in values, so the new "Date" cell
# (G7) ends up sharing the same cell style as the other date cells.
$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A7").Value = 9666.3799999999992
$ws.Range("B7").Value = 9592.52
$ws.Range("C7").Value = 107.96
$ws.Range("D7").Value = 108.79
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = 0.77
$ws.Range("G7").Value = 42609.488483796296
$ws.Range("H7").Value = $true
